$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the styled header row (A1:B1) formatting to the new C1:D1 cells
# before writing values, so the bordered/centered style carries over.
$ws.Range("A1:B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# Row 1 - first interpolation points (bordered style already present on A1:B1,
# now also applied to C1:D1)
$ws.Range("A1").Value = 16.4701252409078
$ws.Range("B1").Value = 158.113202312715
$ws.Range("C1").Value = 490.303442391053
$ws.Range("D1").Value = 493.584887863234

# Row 2
$ws.Range("A2").Value = -1.50764766090166
$ws.Range("B2").Value = -3.68675380357017
$ws.Range("C2").Value = 4.90357404219722
$ws.Range("D2").Value = 8.18501951437844

# Row 3 (new row)
$ws.Range("A3").Value = -0.073593389607704
$ws.Range("B3").Value = -1.53567239662924
$ws.Range("C3").Value = 5.97911474566768
$ws.Range("D3").Value = 8.36427629829022

# Row 4 (new row)
$ws.Range("A4").Value = 8.06264357307232
$ws.Range("B4").Value = -99.1705159487896
$ws.Range("C4").Value = 396.518488954309
$ws.Range("D4").Value = -512.354889313232
